# Fix: Elective lecture and tutorial scheduling
# Applies classroom/lab reassignments across the Regular_Timetable,
# PreMid_Timetable and PostMid_Timetable sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: apply the shared "ELECTIVE BASKETS" block changes (rows 23-32,
# columns D and E) that are identical across the three timetable sheets.
# ---------------------------------------------------------------------------
function Set-ElectiveBasketSlots {
    param($ws)

    $ws.Range("D23").Value = "Tue 13:00-14:30 [C101], Thu 13:00-14:30 [C101]"
    $ws.Range("E23").Value = "Wed 14:30-15:30 [C101]"

    $ws.Range("D24").Value = "Tue 13:00-14:30 [C102], Thu 13:00-14:30 [C102]"
    $ws.Range("E24").Value = "Wed 14:30-15:30 [C102]"

    $ws.Range("D25").Value = "Tue 13:00-14:30 [C104], Thu 13:00-14:30 [C104]"
    $ws.Range("E25").Value = "Wed 14:30-15:30 [C104]"

    $ws.Range("D26").Value = "Tue 13:00-14:30 [C202], Thu 13:00-14:30 [C202]"
    $ws.Range("E26").Value = "Wed 14:30-15:30 [C202]"

    $ws.Range("D27").Value = "Mon 15:30-17:00 [C101], Wed 15:30-17:00 [C101]"
    $ws.Range("E27").Value = "Thu 14:30-15:30 [C101]"

    $ws.Range("D28").Value = "Mon 15:30-17:00 [C102], Wed 15:30-17:00 [C102]"
    $ws.Range("E28").Value = "Thu 14:30-15:30 [C102]"

    $ws.Range("D29").Value = "Mon 15:30-17:00 [C104], Wed 15:30-17:00 [C104]"
    $ws.Range("E29").Value = "Thu 14:30-15:30 [C104]"

    $ws.Range("D30").Value = "Mon 15:30-17:00 [C202], Wed 15:30-17:00 [C202]"
    $ws.Range("E30").Value = "Thu 14:30-15:30 [C202]"

    $ws.Range("D31").Value = "Mon 15:30-17:00 [C203], Wed 15:30-17:00 [C203]"
    $ws.Range("E31").Value = "Thu 14:30-15:30 [C203]"

    $ws.Range("D32").Value = "Mon 15:30-17:00 [C204], Wed 15:30-17:00 [C204]"
    $ws.Range("E32").Value = "Thu 14:30-15:30 [C204]"
}

# ---------------------------------------------------------------------------
# Regular_Timetable
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Regular_Timetable")

$ws.Range("B3").Value = "EC306 [C104]"
$ws.Range("C3").Value = "EC306 [C104]"
$ws.Range("D3").Value = "EC351 [C104]"

$ws.Range("B4").Value = "EC301 [C104]"
$ws.Range("C4").Value = "EC301 [C104]"

$ws.Range("B6").Value = "EC306 (Lab) [L105]"

$ws.Range("B7").Value = "EC306 (Lab) [L105]"
$ws.Range("C7").Value = "EC351 (Tutorial) [C304]"

$ws.Range("B9").Value = "EC301 (Tutorial) [C104]"

Set-ElectiveBasketSlots $ws

# ---------------------------------------------------------------------------
# PreMid_Timetable
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PreMid_Timetable")

$ws.Range("B3").Value = "EC306 [C101]"
$ws.Range("C3").Value = "EC306 [C101]"
$ws.Range("D3").Value = "EC351 [C101]"

$ws.Range("B4").Value = "EC301 [C101]"
$ws.Range("C4").Value = "EC301 [C101]"

$ws.Range("C6").Value = "EC351 [C203]"

$ws.Range("B7").Value = "EC301 (Tutorial) [C104]"

$ws.Range("B8").Value = "EC306 (Lab) [L105]"

$ws.Range("B9").Value = "EC306 (Lab) [L105]"
$ws.Range("C9").Value = "EC351 (Tutorial) [C304]"

Set-ElectiveBasketSlots $ws

# ---------------------------------------------------------------------------
# PostMid_Timetable
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PostMid_Timetable")

$ws.Range("B3").Value = "EC306 [C001]"
$ws.Range("C3").Value = "EC306 [C001]"
$ws.Range("D3").Value = "EC351 [C001]"

$ws.Range("B4").Value = "EC301 [C001]"
$ws.Range("C4").Value = "EC301 [C001]"

$ws.Range("C6").Value = "EC351 [C303]"

$ws.Range("B7").Value = "EC301 (Tutorial) [C101]"

$ws.Range("B8").Value = "EC306 (Lab) [L206]"

$ws.Range("B9").Value = "EC306 (Lab) [L206]"
$ws.Range("C9").Value = "EC351 (Tutorial) [C202]"

Set-ElectiveBasketSlots $ws

$wb.Save()
